$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows (2-13) in this sheet got reshuffled (re-sorted by date) by
# the author. Capture a full snapshot of the columns that vary per-row
# (D, L, M, N, O, P, Q, R, S, T) BEFORE writing anything back, then write
# the rows out in their new order so reads never see already-mutated data.

$srcRows = 2..13
$snapshot = @{}

foreach ($r in $srcRows) {
    $row = @{
        D = $ws.Cells.Item($r, 4).Value2
        L = $ws.Cells.Item($r, 12).Value2
        M = $ws.Cells.Item($r, 13).Value2
        N = $ws.Cells.Item($r, 14).Value2
        O = $ws.Cells.Item($r, 15).Value2
        P = $ws.Cells.Item($r, 16).Value2
        Q = $ws.Cells.Item($r, 17).Value2
        R = $ws.Cells.Item($r, 18).Value2
        S = $ws.Cells.Item($r, 19).Value2
        T = $ws.Cells.Item($r, 20).Value2
    }
    $snapshot[$r] = $row
}

# New row number -> old row number it now holds the data of.
$mapping = @{
    2  = 7
    3  = 4
    4  = 12
    5  = 13
    6  = 9
    7  = 10
    8  = 11
    9  = 2
    10 = 3
    11 = 5
    12 = 6
    13 = 8
}

foreach ($newRow in ($mapping.Keys | Sort-Object)) {
    $oldRow = $mapping[$newRow]
    $data = $snapshot[$oldRow]

    $ws.Cells.Item($newRow, 4).Value = $data.D
    $ws.Cells.Item($newRow, 12).Value = $data.L
    $ws.Cells.Item($newRow, 13).Value = $data.M
    $ws.Cells.Item($newRow, 14).Value = $data.N
    $ws.Cells.Item($newRow, 15).Value = $data.O
    $ws.Cells.Item($newRow, 16).Value = $data.P
    $ws.Cells.Item($newRow, 17).Value = $data.Q
    $ws.Cells.Item($newRow, 18).Value = $data.R
    $ws.Cells.Item($newRow, 19).Value = $data.S
    $ws.Cells.Item($newRow, 20).Value = $data.T
}
